# Update the "dSF" (column F) values for a handful of rows to reflect
# the repulled/recalculated data, per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -8
    4  = -3
    5  = -4
    8  = -1
    9  = -2
    11 = 3
    12 = -5
    13 = -9
    14 = -8
    15 = -5
    16 = -4
    19 = -1
    20 = -6
    21 = 2
    28 = -3
    31 = -4
    34 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
